$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "players" worksheet right before "s_scavenge_vs_gather"
#    (i.e. right after "tile_info"). It mirrors the player roster that
#    already lives on the "ui" sheet (rows 15-37) through live formulas,
#    same pattern the game uses elsewhere to read data off of "ui".
# ---------------------------------------------------------------------
$tileInfo = $wb.Worksheets.Item("tile_info")
$players = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tileInfo)
$players.Name = "players"

$playersRange = $players.Range("A1:N23")
$playersRange.FormulaR1C1 = '=IF(ISBLANK(ui!R[14]C),"",ui!R[14]C)'

$players.Activate()
$players.Range("O1").Select()

# ---------------------------------------------------------------------
# 2. Update view state on a couple of other sheets.
# ---------------------------------------------------------------------

# ui: scroll down a bit (selection itself stays on R5)
$ui = $wb.Worksheets.Item("ui")
$ui.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ui.Range("R5").Select()

# map_tiles: selection moves from P1:Q10 to a single cell G16
$mapTiles = $wb.Worksheets.Item("map_tiles")
$mapTiles.Activate()
$mapTiles.Range("G16").Select()

# ---------------------------------------------------------------------
# 3. Force a full recalculation so volatile formulas (RANDBETWEEN, etc.)
#    refresh their cached values.
# ---------------------------------------------------------------------
$excel.CalculateFullRebuild()

# ---------------------------------------------------------------------
# 4. Leave "s_scavenge_vs_gather" as the active/selected sheet, as in
#    the target workbook (it keeps its own selection at I15).
# ---------------------------------------------------------------------
$scavenge = $wb.Worksheets.Item("s_scavenge_vs_gather")
$scavenge.Activate()
$scavenge.Range("I15").Select()
